$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value2 = 3199.6667
$ws.Range("I64").Value2 = 2899.3333
$ws.Range("J64").Value2 = 3500
$ws.Range("K64").Value2 = 2899.3333
$ws.Range("L64").Value2 = 3500
$ws.Range("M64").Value2 = -2651.3333
$ws.Range("N64").Value2 = -3996
$ws.Range("H67").Value2 = 3199.6667
$ws.Range("I67").Value2 = 2899.3333
$ws.Range("J67").Value2 = 3500
$ws.Range("K67").Value2 = 2899.3333
$ws.Range("L67").Value2 = 3500
$ws.Range("M67").Value2 = -2041.3333
$ws.Range("N67").Value2 = -5216
$ws.Range("H74").Value2 = 5721.7856
$ws.Range("I74").Value2 = 3181.5
$ws.Range("J74").Value2 = 7627
$ws.Range("K74").Value2 = 3181.5
$ws.Range("L74").Value2 = 7627
$ws.Range("M74").Value2 = -2245.5
$ws.Range("N74").Value2 = -9499
$ws.Range("H76").Value2 = 11117396
$ws.Range("I76").Value2 = 7573.5454
$ws.Range("K76").Value2 = 7573.5454
$ws.Range("M76").Value2 = -7258.5454
$ws.Range("H77").Value2 = 5721.7856
$ws.Range("I77").Value2 = 3181.5
$ws.Range("J77").Value2 = 7627
$ws.Range("K77").Value2 = 15907.5
$ws.Range("L77").Value2 = 38135
$ws.Range("M77").Value2 = -11227.5
$ws.Range("N77").Value2 = -47495
$ws.Range("H79").Value2 = 11117396
$ws.Range("I79").Value2 = 7573.5454
$ws.Range("K79").Value2 = 7573.5454
$ws.Range("M79").Value2 = -6481.5454
$ws.Range("H92").Value2 = 1852396.2
$ws.Range("I92").Value2 = 2315133
$ws.Range("J92").Value2 = 1449.1666
$ws.Range("K92").Value2 = 2315133
$ws.Range("L92").Value2 = 1449.1666
$ws.Range("M92").Value2 = -2313885
$ws.Range("N92").Value2 = -3945.1666
$ws.Range("H100").Value2 = 13245
$ws.Range("J100").Value2 = 2855.5557
$ws.Range("L100").Value2 = 2855.5557
$ws.Range("N100").Value2 = -3937.5557
$ws.Range("H113").Value2 = 3313.889
$ws.Range("I113").Value2 = 2891.6667
$ws.Range("K113").Value2 = 2891.6667
$ws.Range("M113").Value2 = 362.3332999999998

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value2 = 18870834
$ws.Range("I132").Value2 = 23810314
$ws.Range("J132").Value2 = 10998.909
$ws.Range("K132").Value2 = 71430942
$ws.Range("L132").Value2 = 32996.727
$ws.Range("M132").Value2 = -71428412
$ws.Range("N132").Value2 = -38056.727

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 1226611.6
$ws.Range("I86").Value2 = 3812.625
$ws.Range("J86").Value2 = 2115920
$ws.Range("K86").Value2 = 3812.625
$ws.Range("L86").Value2 = 2115920
$ws.Range("M86").Value2 = -2689.625
$ws.Range("N86").Value2 = -2118166
$ws.Range("H89").Value2 = 1226611.6
$ws.Range("I89").Value2 = 3812.625
$ws.Range("J89").Value2 = 2115920
$ws.Range("K89").Value2 = 19063.125
$ws.Range("L89").Value2 = 10579600
$ws.Range("M89").Value2 = -13447.125
$ws.Range("N89").Value2 = -10590832
$ws.Range("H94").Value2 = 567.913
$ws.Range("I94").Value2 = 493.44446
$ws.Range("J94").Value2 = 836
$ws.Range("K94").Value2 = 493.44446
$ws.Range("L94").Value2 = 836
$ws.Range("M94").Value2 = -42.44445999999999
$ws.Range("N94").Value2 = -1738
$ws.Range("H105").Value2 = 23811498
$ws.Range("I105").Value2 = 2045.6129
$ws.Range("J105").Value2 = 90910860
$ws.Range("K105").Value2 = 2045.6129
$ws.Range("L105").Value2 = 90910860
$ws.Range("M105").Value2 = -298.6129000000001
$ws.Range("N105").Value2 = -90914354

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value2 = 4091.125
$ws.Range("I62").Value2 = 2895.8
$ws.Range("J62").Value2 = 4634.4546
$ws.Range("K62").Value2 = 2895.8
$ws.Range("L62").Value2 = 4634.4546
$ws.Range("M62").Value2 = -2271.8
$ws.Range("N62").Value2 = -5882.4546
$ws.Range("H65").Value2 = 4091.125
$ws.Range("I65").Value2 = 2895.8
$ws.Range("J65").Value2 = 4634.4546
$ws.Range("K65").Value2 = 14479
$ws.Range("L65").Value2 = 23172.273
$ws.Range("M65").Value2 = -11359
$ws.Range("N65").Value2 = -29412.273
$ws.Range("H105").Value2 = 2309.047
$ws.Range("I105").Value2 = 2248.1584
$ws.Range("J105").Value2 = 3973.3333
$ws.Range("K105").Value2 = 2248.1584
$ws.Range("L105").Value2 = 3973.3333
$ws.Range("M105").Value2 = -501.1583999999998
$ws.Range("N105").Value2 = -7467.3333

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 739.9400000000001
$ws.Range("I131").Value2 = 488.8889
$ws.Range("J131").Value2 = 764.7692
$ws.Range("K131").Value2 = 1466.6667
$ws.Range("L131").Value2 = 2294.3076
$ws.Range("M131").Value2 = 3573.3333
$ws.Range("N131").Value2 = -12374.3076

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 3863.1086
$ws.Range("I70").Value2 = 3866.9644
$ws.Range("J70").Value2 = 3857.111
$ws.Range("K70").Value2 = 3866.9644
$ws.Range("L70").Value2 = 3857.111
$ws.Range("M70").Value2 = -3596.9644
$ws.Range("N70").Value2 = -4397.111
$ws.Range("H73").Value2 = 3863.1086
$ws.Range("I73").Value2 = 3866.9644
$ws.Range("J73").Value2 = 3857.111
$ws.Range("K73").Value2 = 3866.9644
$ws.Range("L73").Value2 = 3857.111
$ws.Range("M73").Value2 = -2930.9644
$ws.Range("N73").Value2 = -5729.111
$ws.Range("H97").Value2 = 552.03845
$ws.Range("I97").Value2 = 518.4167
$ws.Range("J97").Value2 = 955.5
$ws.Range("K97").Value2 = 518.4167
$ws.Range("L97").Value2 = 955.5
$ws.Range("M97").Value2 = -22.41669999999999
$ws.Range("N97").Value2 = -1947.5
$ws.Range("H126").Value2 = 4192.857
$ws.Range("I126").Value2 = 3533.3333
$ws.Range("J126").Value2 = 4687.5
$ws.Range("K126").Value2 = 10599.9999
$ws.Range("L126").Value2 = 14062.5
$ws.Range("M126").Value2 = -8129.999899999999
$ws.Range("N126").Value2 = -19002.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 1988.9615
$ws.Range("I7").Value2 = 1771.8096
$ws.Range("K7").Value2 = 1771.8096
$ws.Range("M7").Value2 = -1659.8096
$ws.Range("H93").Value2 = 1504.2
$ws.Range("I93").Value2 = 1492.75
$ws.Range("J93").Value2 = 1550
$ws.Range("K93").Value2 = 1492.75
$ws.Range("L93").Value2 = 1550
$ws.Range("M93").Value2 = -244.75
$ws.Range("N93").Value2 = -4046
$ws.Range("H126").Value2 = 1988.9615
$ws.Range("I126").Value2 = 1771.8096
$ws.Range("K126").Value2 = 5315.4288
$ws.Range("M126").Value2 = -2845.4288

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value2 = 26325
$ws.Range("J70").Value2 = 28371.428
$ws.Range("L70").Value2 = 28371.428
$ws.Range("N70").Value2 = -29001.428
$ws.Range("H73").Value2 = 26325
$ws.Range("J73").Value2 = 28371.428
$ws.Range("L73").Value2 = 28371.428
$ws.Range("N73").Value2 = -30555.428
$ws.Range("H96").Value2 = 3200
